$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 14) by copying the formatting of the
#     current last data row (row 13) so the new cells pick up the same
#     cell styles (s=4 for the date column, s=3 for the percentage
#     columns) without minting any new style entries. ---
$ws.Range("A13:D13").Copy($ws.Range("A14:D14"))

# --- Update row 13 (previously the last row of data) with its new
#     values from the second data refresh. ---
$ws.Range("A13").Value = 45657
$ws.Range("B13").Value = -0.040456817157672728
$ws.Range("C13").Value = 0.11504492669835839
$ws.Range("D13").Value = 0.057729375378128042

# --- Populate the freshly-created row 14 with the latest data point. ---
$ws.Range("A14").Value = 45658
$ws.Range("B14").Value = 0.033412105974565356
$ws.Range("C14").Value = 0.43513334946172311
$ws.Range("D14").Value = 0.2149864263526432

# --- Move the active selection, matching where the author left off
#     after entering the new row of data. ---
$ws.Range("F16").Select()
